$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("D8").Value = 80.33
$ws.Range("E8").Value = 78.33
$ws.Range("F8").Value = 79.5
$ws.Range("G8").Value = 80.83
$ws.Range("D9").Value = 78.05
$ws.Range("E9").Value = 65.849999999999994
$ws.Range("F9").Value = 67.8
$ws.Range("G9").Value = 69.27
$ws.Range("D17").Value = 79.67
$ws.Range("E17").Value = 78.5
$ws.Range("F17").Value = 78.33
$ws.Range("G17").Value = 80.67
$ws.Range("D18").Value = 78.05
$ws.Range("E18").Value = 65.849999999999994
$ws.Range("F18").Value = 65.849999999999994
$ws.Range("G18").Value = 68.78
$ws.Range("D27").Value = 94.54
$ws.Range("E27").Value = 94.73
$ws.Range("F27").Value = 94.73
$ws.Range("G27").Value = 94.73
$ws.Range("D28").Value = 94.87
$ws.Range("E28").Value = 94.39
$ws.Range("F28").Value = 93.8
$ws.Range("G28").Value = 93.91
$ws.Range("D36").Value = 94.57
$ws.Range("E36").Value = 94.72
$ws.Range("F36").Value = 94.73
$ws.Range("G36").Value = 94.73
$ws.Range("D37").Value = 94.83
$ws.Range("E37").Value = 94.25
$ws.Range("F37").Value = 94.28
$ws.Range("G37").Value = 94.13
$ws.Range("D47").Value = 96.37
$ws.Range("E47").Value = 97.5
$ws.Range("F47").Value = 98.26
$ws.Range("G47").Value = 98.99
$ws.Range("D48").Value = 96.54
$ws.Range("E48").Value = 97.06
$ws.Range("F48").Value = 97.36
$ws.Range("G48").Value = 97.51
$ws.Range("D56").Value = 96.45
$ws.Range("E56").Value = 97.74
$ws.Range("F56").Value = 98.21
$ws.Range("G56").Value = 99.06
$ws.Range("D57").Value = 96.85
$ws.Range("E57").Value = 97.14
$ws.Range("F57").Value = 97.33
$ws.Range("G57").Value = 96.86
$ws.Range("D93").Value = 3.6211E-2
$ws.Range("E93").Value = 1
$ws.Range("F93").Value = 3.3000000000000003E-5
$ws.Range("G93").Value = 1.8E-5
$ws.Range("D94").Value = 3.3066999999999999E-2
$ws.Range("E94").Value = 4.3270000000000003E-2
$ws.Range("F94").Value = 3.9870000000000003E-2
$ws.Range("G94").Value = 3.5126999999999999E-2
$ws.Range("E95").Value = 100
$ws.Range("F95").Value = 100
$ws.Range("G95").Value = 100
$ws.Range("E96").Value = 86.2
$ws.Range("F96").Value = 86.73
$ws.Range("G96").Value = 88
$ws.Range("D102").Value = 4.2464000000000002E-2
$ws.Range("E102").Value = 1.1048000000000001E-2
$ws.Range("F102").Value = 8.8999999999999995E-4
$ws.Range("G102").Value = 6.2000000000000003E-5
$ws.Range("H102").Value = 2.0999999999999999E-5
$ws.Range("I102").Value = 1.8E-5
$ws.Range("D103").Value = 3.9845600000000002E-2
$ws.Range("E103").Value = 2.3820999999999998E-2
$ws.Range("F103").Value = 3.0034999999999999E-2
$ws.Range("G103").Value = 3.3804500000000001E-2
$ws.Range("H103").Value = 0.35009000000000001
$ws.Range("I103").Value = 3.5126999999999999E-2
$ws.Range("D109").Value = 3.5389999999999998E-2
$ws.Range("E109").Value = 1.818055E-2
$ws.Range("F109").Value = 6.3610000000000003E-3
$ws.Range("G109").Value = 8.4099999999999995E-4
$ws.Range("H109").Value = 4.3000000000000002E-5
$ws.Range("I109").Value = 1.9999999999999999E-6
$ws.Range("D110").Value = 3.4311000000000001E-2
$ws.Range("E110").Value = 2.4514999999999999E-2
$ws.Range("F110").Value = 2.4888E-2
$ws.Range("G110").Value = 3.0487E-2
$ws.Range("H110").Value = 3.4533000000000001E-2
$ws.Range("I110").Value = 3.6739000000000001E-2
$ws.Range("F111").Value = 98.79
$ws.Range("E112").Value = 91.6
$ws.Range("H112").Value = 88
$ws.Range("I112").Value = 87.6
